# Generate Report for Handoff
#
# A new handoff was generated for the "85ee349b-2501-4ac1-82f8-284c63e1198b"
# file (row 4 on each locale sheet), which refreshes its "Latest Handoff
# Datetime" (column D) on both the zh-cn and de-de status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D4").Value = "2016-02-26 05:23:13"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D4").Value = "2016-02-26 05:23:25"
